$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-01 Monday", "2025-09-02 Tuesday"),
    @("556÷7=", "702÷6="),
    @("182÷7=", "332÷2="),
    @("504÷9=", "924÷4="),
    @("775÷9=", "957÷8="),
    @("511÷3=", "230÷6="),
    @("749÷5=", "529÷2="),
    @("988÷5=", "685÷4="),
    @("128÷9=", "352÷8="),
    @("390÷7=", "370÷4="),
    @("131÷9=", "238÷3="),
    @("459÷7=", "746÷2="),
    @("974÷6=", "257÷8="),
    @("621÷2=", "107÷6="),
    @("275÷8=", "365÷6="),
    @("167÷8=", "318÷7="),
    @("798÷7=", "985÷2="),
    @("922÷8=", "573÷4="),
    @("886÷2=", "464÷5="),
    @("894÷9=", "521÷8="),
    @("436÷2=", "634÷8="),
    @("374÷3=", "228÷6="),
    @("931÷3=", "354÷2="),
    @("541÷3=", "573÷8="),
    @("158÷8=", "399÷2="),
    @("203÷8=", "782÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
